# "Update 2p0. Convention change to support multi-axle vehicles"
#
# Adds two new worksheets - Truck_Amandla_A2 and Truck_Amandla_A3 - to the
# Differential/Gear1DShafts1D workbook, following the same convention as the
# existing vehicle sheets (e.g. Bus_Makhulu_r). Each new sheet is produced by
# copying the last existing sheet (preserving all data, styles and
# conditional formatting) and then updating the "Instance" name in cell H3 to
# match the new sheet.

$wb = $excel.ActiveWorkbook

# --- Truck_Amandla_A2 : copy of the last sheet (Bus_Makhulu_r) -------------
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$lastSheet.Copy([System.Reflection.Missing]::Value, $lastSheet)
$sheetA2 = $wb.Worksheets.Item($wb.Worksheets.Count)
$sheetA2.Name = "Truck_Amandla_A2"
$sheetA2.Range("H3").Value = "Gear1DShafts1D_Truck_Amandla_A2"
$sheetA2.Activate() | Out-Null
$sheetA2.Range("J9").Select() | Out-Null

# --- Truck_Amandla_A3 : copy of Truck_Amandla_A2 ----------------------------
$sheetA2.Copy([System.Reflection.Missing]::Value, $sheetA2)
$sheetA3 = $wb.Worksheets.Item($wb.Worksheets.Count)
$sheetA3.Name = "Truck_Amandla_A3"
$sheetA3.Range("H3").Value = "Gear1DShafts1D_Truck_Amandla_A3"
$sheetA3.Activate() | Out-Null
$sheetA3.Range("E13").Select() | Out-Null
